# Update "想去人数" (want-to-go count) values in column F
# for the "展览" sheet and the combined "全部类型" sheet.
# These mirror each other (全部类型 aggregates rows from the
# other category sheets), so the same F-column values are bumped
# in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# row -> new value for "展览" sheet
$exhibitUpdates = @{
    4  = 3378
    5  = 222
    6  = 4873
    8  = 308
    12 = 47
    13 = 18
    14 = 668
    20 = 4781
    24 = 5915
    27 = 251
    29 = 4428
    31 = 100
    33 = 883
    35 = 12
    36 = 803
    37 = 876
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# row -> new value for "全部类型" sheet
$allUpdates = @{
    8  = 3378
    9  = 222
    10 = 4873
    12 = 308
    16 = 47
    17 = 18
    18 = 668
    25 = 4781
    29 = 5915
    32 = 251
    34 = 4428
    37 = 100
    39 = 883
    41 = 12
    42 = 803
    43 = 876
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
